$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 10).Value = 2.25
$ws.Cells.Item(2, 11).Value = 2.25
$ws.Cells.Item(2, 13).Value = 1.05
$ws.Cells.Item(2, 14).Value = 11
$ws.Cells.Item(2, 15).Value = 1.29
$ws.Cells.Item(2, 16).Value = 3.75
$ws.Cells.Item(2, 17).Value = 1.92
$ws.Cells.Item(2, 18).Value = 1.98
$ws.Cells.Item(2, 21).Value = 1.91
$ws.Cells.Item(2, 22).Value = 1.91
$ws.Cells.Item(3, 7).Value = 1.91
$ws.Cells.Item(3, 8).Value = 3.2
$ws.Cells.Item(3, 10).Value = 2.75
$ws.Cells.Item(3, 11).Value = 1.91
$ws.Cells.Item(3, 13).Value = 1.13
$ws.Cells.Item(3, 14).Value = 6
$ws.Cells.Item(3, 15).Value = 1.57
$ws.Cells.Item(3, 16).Value = 2.38
$ws.Cells.Item(3, 17).Value = 2.7
$ws.Cells.Item(3, 18).Value = 1.44
$ws.Cells.Item(3, 21).Value = 2.38
$ws.Cells.Item(3, 22).Value = 1.53
$ws.Cells.Item(3, 27).Value = 21
$ws.Cells.Item(3, 29).Value = 6
$ws.Cells.Item(3, 31).Value = 23
$ws.Cells.Item(3, 32).Value = 101
$ws.Cells.Item(3, 34).Value = 8.5
$ws.Cells.Item(3, 42).Value = 29
$ws.Cells.Item(3, 44).Value = 81
$ws.Cells.Item(3, 45).Value = 301
$ws.Cells.Item(3, 46).Value = 2.2
$ws.Cells.Item(3, 47).Value = 10
$ws.Cells.Item(3, 52).Value = 126
$ws.Cells.Item(3, 54).Value = 501
$ws.Cells.Item(4, 7).Value = 1.75
$ws.Cells.Item(4, 9).Value = 5.75
$ws.Cells.Item(4, 10).Value = 2.5
$ws.Cells.Item(4, 13).Value = 1.13
$ws.Cells.Item(4, 14).Value = 6
$ws.Cells.Item(4, 15).Value = 1.5
$ws.Cells.Item(4, 16).Value = 2.5
$ws.Cells.Item(4, 17).Value = 2.6
$ws.Cells.Item(4, 18).Value = 1.48
$ws.Cells.Item(4, 42).Value = 26
$ws.Cells.Item(4, 43).Value = 34
$ws.Cells.Item(5, 26).Value = 29
$ws.Cells.Item(5, 44).Value = 67
$ws.Cells.Item(5, 51).Value = 23
$ws.Cells.Item(5, 53).Value = 67
$ws.Cells.Item(6, 9).Value = 7.4
$ws.Cells.Item(6, 11).Value = 2.5
$ws.Cells.Item(6, 12).Value = 6.5
$ws.Cells.Item(6, 18).Value = 2.22
$ws.Cells.Item(6, 23).Value = 6.9
$ws.Cells.Item(6, 29).Value = 14.5
$ws.Cells.Item(6, 33).Value = 450
$ws.Cells.Item(6, 37).Value = 120
$ws.Cells.Item(6, 38).Value = 60
$ws.Cells.Item(6, 39).Value = 50
$ws.Cells.Item(6, 40).Value = 3.25
$ws.Cells.Item(6, 49).Value = 8.75
$ws.Cells.Item(7, 8).Value = 2.82
$ws.Cells.Item(7, 10).Value = 3.45
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 3
$ws.Cells.Item(7, 14).Value = 6.95
$ws.Cells.Item(7, 15).Value = 1.37
$ws.Cells.Item(7, 16).Value = 2.65
$ws.Cells.Item(7, 17).Value = 2.07
$ws.Cells.Item(7, 18).Value = 1.6
$ws.Cells.Item(7, 21).Value = 1.75
$ws.Cells.Item(7, 22).Value = 1.87
$ws.Cells.Item(7, 23).Value = 8.25
$ws.Cells.Item(7, 24).Value = 15
$ws.Cells.Item(7, 29).Value = 7.5
$ws.Cells.Item(7, 31).Value = 13.5
$ws.Cells.Item(7, 32).Value = 65
$ws.Cells.Item(7, 34).Value = 7.4
$ws.Cells.Item(7, 35).Value = 12.5
$ws.Cells.Item(7, 36).Value = 9.25
$ws.Cells.Item(7, 38).Value = 22
$ws.Cells.Item(7, 39).Value = 32
$ws.Cells.Item(7, 40).Value = 4.9
$ws.Cells.Item(7, 41).Value = 16
$ws.Cells.Item(7, 42).Value = 21
$ws.Cells.Item(7, 43).Value = 75
$ws.Cells.Item(7, 44).Value = 100
$ws.Cells.Item(7, 46).Value = 2.5
$ws.Cells.Item(7, 47).Value = 6.3
$ws.Cells.Item(7, 48).Value = 50
$ws.Cells.Item(7, 49).Value = 4.4
$ws.Cells.Item(7, 50).Value = 13
$ws.Cells.Item(7, 51).Value = 19
$ws.Cells.Item(7, 52).Value = 55
$ws.Cells.Item(7, 53).Value = 80
$ws.Cells.Item(7, 54).Value = 200
$ws.Cells.Item(8, 7).Value = 1.85
$ws.Cells.Item(8, 8).Value = 3.3
$ws.Cells.Item(8, 9).Value = 4.2
$ws.Cells.Item(8, 10).Value = 2.5
$ws.Cells.Item(8, 12).Value = 4.33
$ws.Cells.Item(8, 15).Value = 1.29
$ws.Cells.Item(8, 16).Value = 3.5
$ws.Cells.Item(8, 17).Value = 1.95
$ws.Cells.Item(8, 18).Value = 1.9
$ws.Cells.Item(8, 19).Value = 1.4
$ws.Cells.Item(8, 20).Value = 2.75
$ws.Cells.Item(8, 21).Value = 1.73
$ws.Cells.Item(8, 22).Value = 2
$ws.Cells.Item(8, 23).Value = 7.5
$ws.Cells.Item(8, 24).Value = 9
$ws.Cells.Item(8, 27).Value = 15
$ws.Cells.Item(8, 28).Value = 26
$ws.Cells.Item(8, 29).Value = 10
$ws.Cells.Item(8, 30).Value = 6.5
$ws.Cells.Item(8, 31).Value = 15
$ws.Cells.Item(8, 33).Value = 201
$ws.Cells.Item(8, 34).Value = 13
$ws.Cells.Item(8, 35).Value = 21
$ws.Cells.Item(8, 39).Value = 41
$ws.Cells.Item(8, 41).Value = 10
$ws.Cells.Item(8, 42).Value = 21
$ws.Cells.Item(8, 43).Value = 34
$ws.Cells.Item(8, 45).Value = 151
$ws.Cells.Item(8, 46).Value = 2.75
$ws.Cells.Item(8, 47).Value = 8
$ws.Cells.Item(8, 51).Value = 29
$ws.Cells.Item(8, 54).Value = 201
$ws.Cells.Item(8, 55).Value = 126
$ws.Cells.Item(9, 7).Value = 1.45
$ws.Cells.Item(9, 9).Value = 6.5
$ws.Cells.Item(9, 10).Value = 1.95
$ws.Cells.Item(9, 11).Value = 2.6
$ws.Cells.Item(9, 12).Value = 5.5
$ws.Cells.Item(9, 13).Value = 1.03
$ws.Cells.Item(9, 14).Value = 17
$ws.Cells.Item(9, 15).Value = 1.17
$ws.Cells.Item(9, 16).Value = 5
$ws.Cells.Item(9, 17).Value = 1.53
$ws.Cells.Item(9, 18).Value = 2.4
$ws.Cells.Item(9, 19).Value = 1.25
$ws.Cells.Item(9, 20).Value = 3.75
$ws.Cells.Item(9, 23).Value = 9.5
$ws.Cells.Item(9, 43).Value = 19
$ws.Cells.Item(9, 46).Value = 3.75
$ws.Cells.Item(9, 49).Value = 8
$ws.Cells.Item(9, 52).Value = 101
$ws.Cells.Item(9, 56).Value = 176
$ws.Cells.Item(10, 7).Value = 2.2
$ws.Cells.Item(10, 13).Value = 1.08
$ws.Cells.Item(10, 14).Value = 7.5
$ws.Cells.Item(10, 15).Value = 1.44
$ws.Cells.Item(10, 16).Value = 2.63
$ws.Cells.Item(10, 19).Value = 1.5
$ws.Cells.Item(10, 20).Value = 2.5
$ws.Cells.Item(10, 21).Value = 2
$ws.Cells.Item(10, 22).Value = 1.73
$ws.Cells.Item(10, 23).Value = 6.5
$ws.Cells.Item(10, 24).Value = 9.5
$ws.Cells.Item(10, 25).Value = 9.5
$ws.Cells.Item(10, 27).Value = 21
$ws.Cells.Item(10, 28).Value = 34
$ws.Cells.Item(10, 29).Value = 7.5
$ws.Cells.Item(10, 30).Value = 6
$ws.Cells.Item(10, 31).Value = 17
$ws.Cells.Item(10, 32).Value = 67
$ws.Cells.Item(10, 34).Value = 8.5
$ws.Cells.Item(10, 35).Value = 15
$ws.Cells.Item(10, 36).Value = 13
$ws.Cells.Item(10, 37).Value = 41
$ws.Cells.Item(10, 38).Value = 29
$ws.Cells.Item(10, 39).Value = 41
$ws.Cells.Item(10, 40).Value = 4
$ws.Cells.Item(10, 41).Value = 13
$ws.Cells.Item(10, 42).Value = 26
$ws.Cells.Item(10, 44).Value = 67
$ws.Cells.Item(10, 45).Value = 201
$ws.Cells.Item(10, 46).Value = 2.5
$ws.Cells.Item(10, 47).Value = 8.5
$ws.Cells.Item(10, 48).Value = 67
$ws.Cells.Item(10, 50).Value = 21
$ws.Cells.Item(10, 51).Value = 34
$ws.Cells.Item(10, 52).Value = 67
$ws.Cells.Item(10, 53).Value = 101
$ws.Cells.Item(10, 54).Value = 301
$ws.Cells.Item(11, 13).Value = 1.13
$ws.Cells.Item(11, 14).Value = 6
$ws.Cells.Item(12, 7).Value = 2.1
$ws.Cells.Item(12, 9).Value = 3.5
$ws.Cells.Item(12, 10).Value = 2.88
$ws.Cells.Item(12, 24).Value = 9
$ws.Cells.Item(12, 25).Value = 9.5
$ws.Cells.Item(12, 26).Value = 19
$ws.Cells.Item(12, 28).Value = 34
$ws.Cells.Item(12, 29).Value = 7.5
$ws.Cells.Item(12, 34).Value = 8
$ws.Cells.Item(12, 35).Value = 17
$ws.Cells.Item(12, 36).Value = 13
$ws.Cells.Item(12, 41).Value = 12
$ws.Cells.Item(12, 42).Value = 26
$ws.Cells.Item(12, 44).Value = 67
$ws.Cells.Item(12, 49).Value = 5.5
$ws.Cells.Item(12, 52).Value = 81
$ws.Cells.Item(13, 7).Value = 1.29
$ws.Cells.Item(13, 8).Value = 5
$ws.Cells.Item(13, 9).Value = 12
$ws.Cells.Item(13, 10).Value = 1.73
$ws.Cells.Item(13, 11).Value = 2.5
$ws.Cells.Item(13, 12).Value = 9
$ws.Cells.Item(13, 13).Value = 1.04
$ws.Cells.Item(13, 14).Value = 13
$ws.Cells.Item(13, 15).Value = 1.2
$ws.Cells.Item(13, 16).Value = 4.33
$ws.Cells.Item(13, 17).Value = 1.7
$ws.Cells.Item(13, 18).Value = 2.1
$ws.Cells.Item(13, 19).Value = 1.33
$ws.Cells.Item(13, 20).Value = 3.25
$ws.Cells.Item(13, 23).Value = 6.5
$ws.Cells.Item(13, 25).Value = 9.5
$ws.Cells.Item(13, 26).Value = 7.5
$ws.Cells.Item(13, 28).Value = 34
$ws.Cells.Item(13, 32).Value = 81
$ws.Cells.Item(13, 37).Value = 151
$ws.Cells.Item(13, 38).Value = 81
$ws.Cells.Item(13, 40).Value = 3.2
$ws.Cells.Item(13, 41).Value = 6
$ws.Cells.Item(13, 42).Value = 19
$ws.Cells.Item(13, 43).Value = 15
$ws.Cells.Item(13, 46).Value = 3.25
$ws.Cells.Item(13, 49).Value = 10
$ws.Cells.Item(13, 52).Value = 251
$ws.Cells.Item(13, 53).Value = 251
